# Update localization status report: rows for
#   7e6a22d9-868a-4396-8c40-6fb8d922c1f3  (row 3)
#   91b01dfe-d9d4-4034-bf7c-7b77f8304f58  (row 4)
# move from "Ready for handoff" to "In Translation".
#
# This touches the Status column ("C") on the per-language sheets
# ("zh-cn" and "de-de"), and the corresponding language-status columns
# ("B" = zh-cn, "C" = de-de) on the "Overview" sheet.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: zh-cn (col B) and de-de (col C) status for rows 3 & 4 ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B3").Value = $newStatus
$ov.Range("C3").Value = $newStatus
$ov.Range("B4").Value = $newStatus
$ov.Range("C4").Value = $newStatus

# --- zh-cn sheet: Status column (C) for rows 3 & 4 ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = $newStatus
$zh.Range("C4").Value = $newStatus

# --- de-de sheet: Status column (C) for rows 3 & 4 ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = $newStatus
$de.Range("C4").Value = $newStatus
